# Updates the cryptos price list (Price/Volume columns) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "63.405.22"

$ws.Range("D3").Value = "3.099.28"
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue "D5" "552.91"
$ws.Range("E5").Value = "  -2.59%  "

Set-TextValue "D6" "138.13"
$ws.Range("E6").Value = "  -8.99%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.090.16"
$ws.Range("E8").Value = "  -1.79%  "

Set-TextValue "D9" "0.497"
$ws.Range("E9").Value = "  -1.48%  "

Set-TextValue "D10" "0.161"
$ws.Range("E10").Value = "  -0.85%  "

Set-TextValue "D11" "6.56"
$ws.Range("E11").Value = "  -2.74%  "

Set-TextValue "D12" "0.460"
$ws.Range("E12").Value = "  -1.98%  "

Set-TextValue "D13" "35.10"
$ws.Range("E13").Value = "  -6.35%  "

$ws.Range("E14").Value = "  -4.05%  "

$ws.Range("D15").Value = "3.600.66"

$ws.Range("D16").Value = "63.414.80"
$ws.Range("E16").Value = "  -3.46%  "

$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").Value = "3.100.43"
$ws.Range("E18").Value = "  -1.84%  "

Set-TextValue "D19" "506.79"
$ws.Range("E19").Value = "  -4.19%  "

$ws.Range("E20").Value = "  -2.55%  "

Set-TextValue "D21" "13.55"
$ws.Range("E21").Value = "  -4.07%  "

Set-TextValue "D22" "0.706"
$ws.Range("E22").Value = "  -0.94%  "

Set-TextValue "D23" "7.24"
$ws.Range("E23").Value = "  -3.89%  "

Set-TextValue "D24" "77.98"
$ws.Range("E24").Value = "  -2.46%  "

Set-TextValue "D25" "12.32"
$ws.Range("E25").Value = "  -4.75%  "

$ws.Range("E26").Value = "  +0.11%  "

Set-TextValue "D27" "2.75"
$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("E28").Value = "  -7.83%  "

$ws.Range("E29").Value = "  +0.02%  "

Set-TextValue "D30" "1.95"
$ws.Range("E30").Value = "  -10.67%  "

Set-TextValue "D31" "26.45"
$ws.Range("E31").Value = "  -2.15%  "

Set-TextValue "D32" "2.53"
$ws.Range("E32").Value = "  -8.15%  "

Set-TextValue "D33" "1.12"
$ws.Range("E33").Value = "  -2.75%  "

Set-TextValue "D34" "525.67"
$ws.Range("E34").Value = "  -10.11%  "

Set-TextValue "D35" "57.34"
$ws.Range("E35").Value = "  +7.47%  "

Set-TextValue "D36" "6.00"
$ws.Range("E36").Value = "  -2.83%  "

Set-TextValue "D37" "5.21"
$ws.Range("E37").Value = "  -8.83%  "

Set-TextValue "D38" "0.0414"
$ws.Range("E38").Value = "  -3.58%  "

$ws.Range("D39").Value = "3.080.05"
$ws.Range("E39").Value = "  +0.72%  "

Set-TextValue "D40" "0.0793"
$ws.Range("E40").Value = "  -5.09%  "

Set-TextValue "D41" "0.120"
$ws.Range("E41").Value = "  -2.47%  "

Set-TextValue "D42" "2.75"
$ws.Range("E42").Value = "  -10.01%  "

Set-TextValue "D43" "8.13"
$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D44" "0.253"
$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue "D45" "2.70"
$ws.Range("E45").Value = "  +74.48%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D47" "123.18"
$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D48" "2.05"
$ws.Range("E48").Value = "  -7.75%  "

Set-TextValue "D49" "24.35"
$ws.Range("E49").Value = "  -7.06%  "

Set-TextValue "D50" "0.107"
$ws.Range("E50").Value = "  -2.99%  "

$ws.Range("D51").Value = "0.0₃0509"
$ws.Range("E51").Value = "  -7.59%  "
